$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels per editor comments: "means" -> "mean (n)", "sds" -> "sd (n)"
$ws.Range("C1").Value = "mean (n)"
$ws.Range("D1").Value = "sd (n)"
